$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update existing row 10 (FW_UI_0008 / SetHeadlineBold):
#   - D10 was blank, now holds "ON"
#   - F10 description text changes
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "ON"
$ws.Range("F10").Value = "To verify that user is able to set headline in bold"

# ---------------------------------------------------------------------------
# Update existing row 11 (FW_UI_0009 / EditExistingHeadlineHighlights):
#   - A11 "Yes" -> "No"
#   - F11 description text changes
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "No"
$ws.Range("F11").Value = "To verify that user is able to set color,font options and preview headline for an existing headline alarm based on the selections made"

# ---------------------------------------------------------------------------
# New row 12 (FW_UI_0010) - same banded style as row 10 (s=2)
# ---------------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A12").Value = "No"
$ws.Range("B12").Value = "FW_UI_0010"
$ws.Range("C12").Value = "HeadlineOptionReopenRelaunch"
$ws.Range("D12").Value = "Reopen"
$ws.Range("E12").Value = "HeadLineSettings22797"
$ws.Range("F12").Value = "To verify that user is able to view headline highlight settings done under alarms list on reopening preferences"

# ---------------------------------------------------------------------------
# New row 13 (FW_UI_0011) - banded style like row 11 (s=5)
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$ws.Range("A13").Value = "No"
$ws.Range("B13").Value = "FW_UI_0011"
$ws.Range("C13").Value = "HeadlineOptionReopenRelaunch"
$ws.Range("D13").Value = "Relaunch"
$ws.Range("E13").Value = "HeadLineSettings22797"
$ws.Range("F13").Value = "To verify that headline highlight settings for new and existing alarms are retained on relaunching LE"

# ---------------------------------------------------------------------------
# New row 14 (FW_UI_0012) - also s=5 (matches source workbook exactly,
# the banding does not alternate back to s=2 here)
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A14").Value = "No"
$ws.Range("B14").Value = "FW_UI_0012"
$ws.Range("C14").Value = "SetHeadlineBold"
$ws.Range("D14").Value = "OFF"
$ws.Range("E14").Value = "HeadLineSettings22797"
$ws.Range("F14").Value = "To verify that user is able to set headline without bold option"

# ---------------------------------------------------------------------------
# New row 15 (FW_UI_0013) - back to s=2
# ---------------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Yes"
$ws.Range("B15").Value = "FW_UI_0013"
$ws.Range("C15").Value = "VerifyFeedsDropdown"
$ws.Range("E15").Value = "Feeds233961"
$ws.Range("F15").Value = "To verify that user is able to view the the selected feeds under Feeds dropdown in UI"

# ---------------------------------------------------------------------------
# New blank rows 16-24, continuing the alternating s=5 / s=2 band styling
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A21:F21").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A23:F23").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Stamp an explicit row height on every newly added row so each one carries
# the same ht="13.5" customHeight="1" attributes as the pre-existing rows.
for ($r = 12; $r -le 24; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.5
}

# Selection ends on C15, matching the saved worksheet view.
$ws.Range("C15").Select()
